$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 7 through 22 (columns A-H). This both rewrites the
# existing rows 7-10 with their new values and appends the additional
# spans as rows 11-22, per the "add the other spans" commit.
$rows = @(
    @(-8, 2.5,  6.6259, 0, 250, 112, 0.16, 0.158),
    @(-1, 2.5,  3.4017, 0, 250, 112, 0.16, 0.158),
    @(2,  2.5,  -2.1036, 0, 250, 112, 0.16, 0.158),
    @(-2, 7.5,  6.6443, 0, 250, 112, 0.16, 0.158),
    @(-1, 7.5,  5.6105, 0, 250, 112, 0.16, 0.158),
    @(-8, 10,   6.2596, 0, 250, 112, 0.16, 0.158),
    @(-8, 12.5, 6.0851, 0, 250, 112, 0.16, 0.158),
    @(-5, 12.5, 7.6323, 0, 250, 112, 0.16, 0.158),
    @(-7, 15,   6.461,  0, 250, 112, 0.16, 0.158),
    @(-1, 17.5, 6.9065, 0, 250, 112, 0.16, 0.158),
    @(2,  17.5, 2.4518, 0, 250, 112, 0.16, 0.158),
    @(-8, 20,   5.7182, 0, 250, 112, 0.16, 0.158),
    @(-7, 20,   6.2239, 0, 250, 112, 0.16, 0.158),
    @(-6, 20,   6.6515, 0, 250, 112, 0.16, 0.158),
    @(-6, 25,   6.3423, 0, 250, 112, 0.16, 0.158),
    @(2,  25,   4.5182, 0, 250, 112, 0.16, 0.158)
)

$startRow = 7
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $values = $rows[$i]
    for ($c = 0; $c -lt $values.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $values[$c]
    }
}
